$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column N.
# This shifts the existing mfd_hab1/mfd_hab2/mfd_hab3 columns (N,O,P) to (Q,R,S)
# and leaves three new blank (style-inherited) header cells at N1:P1.
$ws.Range("N1:P1").EntireColumn.Insert()

# Fill in the new header cells
$ws.Range("N1").Value = "correct_mfd_sampletype"
$ws.Range("O1").Value = "correct_mfd_areatype"
$ws.Range("P1").Value = "correct_habitat_typenumber"

# Columns F, M and P hold numeric-looking codes (habitat type numbers / sampling
# date codes) that must be stored as literal text, not auto-converted numbers.
# Format those ranges as Text first, then reset the style back to Normal after
# writing so no visible formatting difference remains.
# (Multi-area ranges only apply NumberFormat/Style to their first area in this
# engine, so each column range is handled individually.)
$ws.Range("F2:F17").NumberFormat = "@"
$ws.Range("M2:M17").NumberFormat = "@"
$ws.Range("P2:P17").NumberFormat = "@"

# --- Rows 2-9: natural_soil / lille vildmose samples ---
$soilRows = 2,3,4,5,6,7,8,9
foreach ($r in $soilRows) {
    $ws.Range("F$r").Value = "7000"
    $ws.Range("K$r").Value = "Soil"
    $ws.Range("L$r").Value = "Natural"
    $ws.Range("M$r").Value = "14169"
    $ws.Range("N$r").Value = "Soil"
    $ws.Range("O$r").Value = "Natural"
    $ws.Range("P$r").Value = "7000"
    $ws.Range("Q$r").Value = "Bogs, mires and fens"
}

# --- Rows 10-13, 15-17: built_environment / thingbæk kalkminer samples ---
$chalkRows = 10,11,12,13,15,16,17
foreach ($r in $chalkRows) {
    $ws.Range("F$r").Value = "2130"
    $ws.Range("K$r").Value = "Sediment"
    $ws.Range("L$r").Value = "Urban"
    $ws.Range("M$r").Value = "14169"
    $ws.Range("N$r").Value = "Sediment"
    $ws.Range("O$r").Value = "Urban"
    $ws.Range("P$r").Value = "2130"
    $ws.Range("Q$r").Value = "Urban"
    $ws.Range("R$r").Value = "Other"
    $ws.Range("S$r").Value = "High chalk concentration (limestone quarry)"
}

# --- Row 14: built_environment / læsø saltsydekar sample ---
$ws.Range("F14").Value = "2120"
$ws.Range("K14").Value = "Sediment"
$ws.Range("L14").Value = "Urban"
$ws.Range("M14").Value = "13984"
$ws.Range("N14").Value = "Sediment"
$ws.Range("O14").Value = "Urban"
$ws.Range("P14").Value = "2120"
$ws.Range("Q14").Value = "Urban"
$ws.Range("R14").Value = "Other"
$ws.Range("S14").Value = "High salinity (saltworks)"

# Restore the default (Normal) style on the text-forced columns so no stray
# number formatting is left behind.
$ws.Range("F2:F17").Style = "Normal"
$ws.Range("M2:M17").Style = "Normal"
$ws.Range("P2:P17").Style = "Normal"
